# Reworking the admin part
$wb = $excel.ActiveWorkbook

# Rename sheets: drop the "Группа " prefix from the 2nd and 3rd sheets
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "А-2-24"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Ю-2-24"

# Make the 3rd sheet ("Ю-2-24") the active tab
$ws3.Activate()
